# "Updated symbol list" refresh (coin prices in column D, a couple of
# relabeled "Worst/Best in 24h" tags in column E, and the "Hora" counter in
# column G bumped from 0 to 2 for every data row). Rows 49/50 (BOLO /
# CryptobidCoin) also swapped places with refreshed price data.
#
# All of these cells are stored as literal TEXT (inlineStr) in the sheet,
# even the ones that look numeric (prices, the "0"/"2" hour flag, etc.), so
# each numeric-looking value is entered with a leading single-quote to force
# Excel to keep it as text instead of silently re-typing the cell as a
# Number. The Style is then reset to "Normal" so the quote-prefix marker
# doesn't leave a stray cell style behind - this mirrors the original file,
# where none of these data cells carry an explicit style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    # Force a literal text entry (handles values that look like numbers)
    # without leaving the quote-prefix cell style attached afterwards.
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

function Set-StringCell($addr, $value) {
    # Plain (non-numeric-looking) text can be assigned directly.
    $ws.Range($addr).Value = $value
}

# Column D (Price) + Column G (Hora) updates for every data row, 2..51.
$priceUpdates = [ordered]@{
    "D2"  = "272.71"
    "D3"  = "23.38"
    "D4"  = "6.402"
    "D5"  = "0.06273"
    "D6"  = "3.669"
    "D7"  = "6.672"
    "D9"  = "0.8355"
    "D10" = "0.01381"
    "D11" = "0.1602"
    "D12" = "0.08322"
    "D13" = "0.03450"
    "D14" = "0.03168"
    "D15" = "0.09301"
    "D16" = "3.839"
    "D17" = "0.001642"
    "D18" = "0.04751"
    "D19" = "0.006368"
    "D20" = "0.005686"
    "D21" = "0.001078"
    "D22" = "0.0001500"
    "D23" = "3.715"
    "D24" = "2.373"
    "D25" = "0.3348"
    "D26" = "0.1275"
    "D27" = "0.0002681"
    "D40" = "0.04725"
    "D41" = "0.007044"
    "D42" = "0.1164"
    "D43" = "0.003700"
    "D44" = "0.01206"
    "D45" = "0.00006266"
    "D47" = "0.00000000750"
    "D48" = "0.7968"
    "D49" = "0.00002300"
    "D50" = "0.002124"
    "D51" = "0.01240"
}

foreach ($addr in $priceUpdates.Keys) {
    Set-TextCell $addr $priceUpdates[$addr]
}

# Column G ("Hora") goes from "0" to "2" on every data row (2..51).
for ($row = 2; $row -le 51; $row++) {
    Set-TextCell "G$row" "2"
}

# Column E label tweaks.
Set-StringCell "E43" "42CEJICEJI"
Set-StringCell "E48" "47CoinbaseStockTokenCOINWorstin24h"
Set-StringCell "E49" "48CryptobidCoinCBC"
Set-StringCell "E50" "49BOLOBOLOBestin24h"

# Rows 49/50 swap identity: BOLO <-> CryptobidCoin (name, link).
Set-StringCell "B49" "CryptobidCoin"
Set-StringCell "C49" "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"

Set-StringCell "B50" "BOLO"
Set-StringCell "C50" "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
